$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily COVID records appended after the last existing row (285),
# matching the layout/columns of the pre-existing data rows:
# data, dia desde 1 contagio, casos, mortes, Ativos, taxa morte contaminados,
# Curados, Casos negativos, Testes realizados, novosCasos, suspeitos,
# mortesSuspeitas, suspeitosAtivos, novosTestes, leitos_clinicos_ocupados,
# leitos_uti_ocupados, novasMortes, semana
$newRows = @(
    @{ Row = 286; Data = @("08/01/2021", 285, 5210, 120, 161, 0.02303262955854127, 4911, 14698, 19908, 85, 490, 1, 489, 193, 25, 10, 4, 41) },
    @{ Row = 287; Data = @("09/01/2021", 286, 5252, 120, 133, 0.02284843869002285, 4981, 14760, 20012, 42, 468, 2, 466, 104, 28, 10, 0, 41) },
    @{ Row = 288; Data = @("10/01/2021", 287, 5265, 120, 108, 0.02279202279202279, 5019, 14760, 20025, 13, 469, 3, 466, 13, 26, 12, 0, 42) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $values = $entry.Data

    # Column A holds the date as literal text (e.g. "08/01/2021"), exactly
    # like the existing rows - force text formatting first so Excel doesn't
    # silently reinterpret it as a date serial number, then drop back to the
    # default "Normal" style so no explicit style sticks to the cell (same
    # as the rest of the data rows, which carry no "s" attribute).
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $values[0]
    $dateCell.Style = "Normal"

    for ($col = 2; $col -le $values.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $values[$col - 1]
    }
}
